{"js": "// Insert a new bulleted list item \"Ola Hansen\" right after the\n// \"Jonas Fredriksen (Eksternt byr\u00e5)\" participant in the\n// \"Protokoll fra oppstartsm\u00f8te\" meeting's \"Deltakere\" list.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the paragraph whose text is the \"Jonas Fredriksen\" list entry.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.trim() === \"Jonas Fredriksen (Eksternt byr\u00e5)\") {\n    target = para;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find paragraph \"Jonas Fredriksen (Eksternt byr\u00e5)\".');\n}\n\n// Insert a new paragraph right after it; Word carries over the list\n// formatting (numId/ilvl) from the preceding list paragraph automatically.\nconst newPara = target.insertParagraph(\"Ola Hansen\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Insert a new bulleted list item \"Ola Hansen\" right after the\n# \"Jonas Fredriksen (Eksternt byr\u00e5)\" participant in the\n# \"Protokoll fra oppstartsm\u00f8te\" meeting's \"Deltakere\" list.\n\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$searchRange.Find.ClearFormatting()\n$found = $searchRange.Find.Execute(\"Jonas Fredriksen (Eksternt byr\u00e5)\")\n\nif (-not $found) {\n    throw 'Could not find paragraph \"Jonas Fredriksen (Eksternt byr\u00e5)\".'\n}\n\n# $searchRange now spans the matched text; grab its enclosing paragraph.\n$targetPara = $searchRange.Paragraphs.Item(1)\n$targetIndex = $targetPara.Range.Paragraphs.Item(1).Index\n\n# Insert a new paragraph right after it; Word carries over the list\n# formatting (numId/ilvl) from the preceding list paragraph automatically.\n$targetPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($targetIndex + 1)\n$newPara.Range.Text = \"Ola Hansen\"\n"}
